$wb = $excel.ActiveWorkbook

# --- Section_A (sheet1) ---
$wsA = $wb.Worksheets.Item("Section_A")

$wsA.Range("B2").Value = "C202"
$wsA.Range("C2").Value = "EC161"
$wsA.Range("D2").Value = "Free"
$wsA.Range("E2").Value = "MA162"

$wsA.Range("B3").Value = "DS161"
$wsA.Range("C3").Value = "MA161"
$wsA.Range("D3").Value = "EC161"
$wsA.Range("E3").Value = "EC161"

$wsA.Range("B5").Value = "EC161"
$wsA.Range("C5").Value = "CS161"
$wsA.Range("D5").Value = "DS161"
$wsA.Range("E5").Value = "Free"
$wsA.Range("F5").Value = "Free"

$wsA.Range("B7").Value = "MA162"
$wsA.Range("C7").Value = "C202"
$wsA.Range("D7").Value = "CS161"
$wsA.Range("E7").Value = "MA161"
$wsA.Range("F7").Value = "CS161"

# --- Section_B (sheet2) ---
$wsB = $wb.Worksheets.Item("Section_B")

$wsB.Range("B2").Value = "CS161"
$wsB.Range("C2").Value = "C202"
$wsB.Range("D2").Value = "CS161"
$wsB.Range("E2").Value = "MA162"

$wsB.Range("B3").Value = "MA161"
$wsB.Range("C3").Value = "DS161"
$wsB.Range("D3").Value = "Free"
$wsB.Range("E3").Value = "C202"

$wsB.Range("B5").Value = "C202"
$wsB.Range("D5").Value = "MA162"
$wsB.Range("E5").Value = "EC161"
$wsB.Range("F5").Value = "Free"

$wsB.Range("B7").Value = "EC161"
$wsB.Range("C7").Value = "MA161"
$wsB.Range("D7").Value = "EC161"
$wsB.Range("E7").Value = "CS161"
$wsB.Range("F7").Value = "DS161"
